$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '41.396.82'
$ws.Range("E2").Value = '  -0.66%  '

# Row 3
$ws.Range("D3").Value = '2.186.36'
$ws.Range("E3").Value = '  -1.36%  '

# Row 4
$ws.Range("E4").Value = '  -0.22%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '249.81'
$ws.Range("E5").Value = '  -1.65%  '

# Row 6
$ws.Range("E6").Value = '  -1.93%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '67.24'
$ws.Range("E7").Value = '  -4.12%  '

# Row 8
$ws.Range("E8").Value = '  -0.12%  '

# Row 9
$ws.Range("E9").Value = '  +3.43%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.65'
$ws.Range("E10").Value = '  -2.09%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '59.48'
$ws.Range("E11").Value = '  +1.36%  '

# Row 12
$ws.Range("E12").Value = '  -3.12%  '

# Row 13
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.104'
$ws.Range("E13").Value = '  -1.16%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.93'
$ws.Range("E14").Value = '  -4.43%  '

# Row 15
$ws.Range("D15").Value = '2.514.77'
$ws.Range("E15").Value = '  -1.29%  '

# Row 16
$ws.Range("E16").Value = '  -3.14%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.856'
$ws.Range("E17").Value = '  -4.63%  '

# Row 18
$ws.Range("D18").Value = '2.176.56'
$ws.Range("E18").Value = '  -2.11%  '

# Row 19
$ws.Range("D19").Value = '41.274.13'
$ws.Range("E19").Value = '  -0.95%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0952'
$ws.Range("E20").Value = '  -1.41%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.12'
$ws.Range("E21").Value = '  -2.60%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.93'
$ws.Range("E22").Value = '  -0.85%  '

# Row 23
$ws.Range("E23").Value = '  -1.61%  '

# Row 24
$ws.Range("E24").Value = '  -1.31%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.86'
$ws.Range("E25").Value = '  -3.52%  '

# Row 26
$ws.Range("E26").Value = '  -0.04%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.33'
$ws.Range("E27").Value = '  -5.60%  '

# Row 28
$ws.Range("E28").Value = '  -5.42%  '

# Row 29
$ws.Range("E29").Value = '  -2.70%  '

# Row 30
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '166.78'
$ws.Range("E30").Value = '  -3.11%  '

# Row 31
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.05'
$ws.Range("E31").Value = '  -6.53%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.23'
$ws.Range("E32").Value = '  -2.58%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0783'
$ws.Range("E33").Value = '  +5.17%  '

# Row 34
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.84'
$ws.Range("E34").Value = '  +3.02%  '

# Row 35
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.120'
$ws.Range("E35").Value = '  -3.55%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.121'
$ws.Range("E36").Value = '  -2.91%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.17'
$ws.Range("E37").Value = '  +3.73%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.79'
$ws.Range("E38").Value = '  -1.22%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.56'
$ws.Range("E39").Value = '  -3.13%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0304'
$ws.Range("E40").Value = '  +0.41%  '

# Row 41
$ws.Range("E41").Value = '  -2.79%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.06'
$ws.Range("E42").Value = '  -2.02%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.14'
$ws.Range("E43").Value = '  +3.60%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.59'
$ws.Range("E44").Value = '  -5.38%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.54'
$ws.Range("E45").Value = '  -5.15%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.195'
$ws.Range("E46").Value = '  -4.15%  '

# Row 47
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1000'
$ws.Range("E47").Value = '  -2.03%  '

# Row 48
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.50'
$ws.Range("E48").Value = '  -3.61%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.999'
$ws.Range("E49").Value = '  -0.29%  '

# Row 50
$ws.Range("E50").Value = '  -1.71%  '

# Row 51
$ws.Range("E51").Value = '  +5.74%  '
